$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$arr = New-Object "object[,]" 24,14

$arr[0,0] = 9.038004860040518
$arr[0,1] = 4.710672694227462
$arr[0,2] = 11.57035918323165
$arr[0,3] = 0
$arr[0,4] = 29.40894198588861
$arr[0,5] = 27.26084037561257
$arr[0,6] = 13.99981090534322
$arr[0,7] = 19.8349590857142
$arr[0,8] = 11.31380303013707
$arr[0,9] = 8.917692391236542
$arr[0,10] = 0
$arr[0,11] = 0
$arr[0,12] = 17.83658848248331
$arr[0,13] = 21.09430479861086
$arr[1,0] = 8.725188728770592
$arr[1,1] = 4.491379053372159
$arr[1,2] = 11.46712753918362
$arr[1,3] = 0
$arr[1,4] = 29.43958170921778
$arr[1,5] = 27.32847748663854
$arr[1,6] = 14.04275012866949
$arr[1,7] = 19.91801167397596
$arr[1,8] = 11.29318466548975
$arr[1,9] = 8.703783438734435
$arr[1,10] = 0
$arr[1,11] = 0
$arr[1,12] = 17.88507821839384
$arr[1,13] = 21.16380702032786
$arr[2,0] = 8.52827660573047
$arr[2,1] = 4.350242072843333
$arr[2,2] = 11.40570169782628
$arr[2,3] = 0
$arr[2,4] = 29.46539221252046
$arr[2,5] = 27.37812135511141
$arr[2,6] = 14.07110157296353
$arr[2,7] = 19.97242733456552
$arr[2,8] = 11.28275117605923
$arr[2,9] = 8.570723280789245
$arr[2,10] = 0
$arr[2,11] = 0
$arr[2,12] = 17.9163726265988
$arr[2,13] = 21.21056037328835
$arr[3,0] = 8.446937889412943
$arr[3,1] = 4.291140142259653
$arr[3,2] = 11.38118778172982
$arr[3,3] = 0
$arr[3,4] = 29.4776681810246
$arr[3,5] = 27.40038487306369
$arr[3,6] = 14.08315478630467
$arr[3,7] = 19.99546255678096
$arr[3,8] = 11.27906245964248
$arr[3,9] = 8.516143588001391
$arr[3,10] = 0
$arr[3,11] = 0
$arr[3,12] = 17.929508914834
$arr[3,13] = 21.230637073155
$arr[4,0] = 8.433369285134566
$arr[4,1] = 4.281231880591119
$arr[4,2] = 11.37714924013112
$arr[4,3] = 0
$arr[4,4] = 29.47981271207971
$arr[4,5] = 27.40420429186729
$arr[4,6] = 14.08518640479001
$arr[4,7] = 19.99933950885855
$arr[4,8] = 11.27848404046743
$arr[4,9] = 8.507061414438063
$arr[4,10] = 0
$arr[4,11] = 0
$arr[4,12] = 17.9317133802563
$arr[4,13] = 21.23403262217348
$arr[5,0] = 8.527183909112019
$arr[5,1] = 4.349451366876517
$arr[5,2] = 11.40536896640976
$arr[5,3] = 0
$arr[5,4] = 29.46555065523084
$arr[5,5] = 27.37841338623118
$arr[5,6] = 14.07126210296477
$arr[5,7] = 19.97273451168591
$arr[5,8] = 11.28269914521813
$arr[5,9] = 8.569988541227994
$arr[5,10] = 0
$arr[5,11] = 0
$arr[5,12] = 17.91654823285951
$arr[5,13] = 21.21082698896133
$arr[6,0] = 8.931220516291786
$arr[6,1] = 4.636431612045912
$arr[6,2] = 11.53437368237424
$arr[6,3] = 0
$arr[6,4] = 29.41805345016414
$arr[6,5] = 27.28247338805135
$arr[6,6] = 14.01420416974938
$arr[6,7] = 19.8628856708916
$arr[6,8] = 11.30623380952493
$arr[6,9] = 8.844337715203816
$arr[6,10] = 0
$arr[6,11] = 0
$arr[6,12] = 17.85299253531311
$arr[6,13] = 21.11742186363863
$arr[7,0] = 9.680323616639797
$arr[7,1] = 5.146106509384078
$arr[7,2] = 11.80175006482308
$arr[7,3] = 0
$arr[7,4] = 29.3804778220624
$arr[7,5] = 27.15901043233234
$arr[7,6] = 13.91806624541577
$arr[7,7] = 19.67460982653501
$arr[7,8] = 11.36989000416492
$arr[7,9] = 9.365557374635019
$arr[7,10] = 0
$arr[7,11] = 0
$arr[7,12] = 17.74038813460367
$arr[7,13] = 20.96667396946135
$arr[8,0] = 10.19869641500464
$arr[8,1] = 5.486479714497356
$arr[8,2] = 12.00538786466502
$arr[8,3] = 0
$arr[8,4] = 29.38676212960695
$arr[8,5] = 27.10808071562344
$arr[8,6] = 13.8570228174223
$arr[8,7] = 19.55281742442365
$arr[8,8] = 11.42708828911151
$arr[8,9] = 9.734356706638438
$arr[8,10] = 0
$arr[8,11] = 0
$arr[8,12] = 17.66492771799204
$arr[8,13] = 20.87575839318865
$arr[9,0] = 10.42658381151412
$arr[9,1] = 5.633660648307116
$arr[9,2] = 12.09924429487851
$arr[9,3] = 0
$arr[9,4] = 29.39696716164611
$arr[9,5] = 27.09360707372123
$arr[9,6] = 13.8313319072193
$arr[9,7] = 19.50099963626508
$arr[9,8] = 11.45531345832883
$arr[9,9] = 9.898322286582616
$arr[9,10] = 0
$arr[9,11] = 0
$arr[9,12] = 17.63216421004669
$arr[9,13] = 20.83872158477329
$arr[10,0] = 10.51167001011511
$arr[10,1] = 5.688276503432147
$arr[10,2] = 12.13493141259931
$arr[10,3] = 0
$arr[10,4] = 29.40188533301499
$arr[10,5] = 27.08937983583055
$arr[10,6] = 13.82190207063473
$arr[10,7] = 19.48189333056801
$arr[10,8] = 11.46631289102999
$arr[10,9] = 9.95980930056426
$arr[10,10] = 0
$arr[10,11] = 0
$arr[10,12] = 17.6199814307864
$arr[10,13] = 20.82531928037377
$arr[11,0] = 10.49340001255241
$arr[11,1] = 5.676563982826444
$arr[11,2] = 12.12723954625023
$arr[11,3] = 0
$arr[11,4] = 29.40077930235692
$arr[11,5] = 27.09023444886882
$arr[11,6] = 13.8239196685584
$arr[11,7] = 19.48598526286844
$arr[11,8] = 11.46393022612007
$arr[11,9] = 9.946594647566673
$arr[11,10] = 0
$arr[11,11] = 0
$arr[11,12] = 17.62259525856867
$arr[11,13] = 20.82817798939439
$arr[12,0] = 10.43360852383859
$arr[12,1] = 5.638176428988966
$arr[12,2] = 12.10217754557575
$arr[12,3] = 0
$arr[12,4] = 29.39735068209345
$arr[12,5] = 27.09323415612172
$arr[12,6] = 13.83055012179674
$arr[12,7] = 19.49941740656489
$arr[12,8] = 11.45621218748861
$arr[12,9] = 9.90339325554336
$arr[12,10] = 0
$arr[12,11] = 0
$arr[12,12] = 17.63115743983942
$arr[12,13] = 20.83760648112972
$arr[13,0] = 10.39682501284405
$arr[13,1] = 5.614516840937557
$arr[13,2] = 12.08684444014778
$arr[13,3] = 0
$arr[13,4] = 29.39538768908642
$arr[13,5] = 27.09523490471616
$arr[13,6] = 13.83465037044135
$arr[13,7] = 19.50771218521674
$arr[13,8] = 11.45152500607808
$arr[13,9] = 9.876851008183166
$arr[13,10] = 0
$arr[13,11] = 0
$arr[13,12] = 17.63643117680582
$arr[13,13] = 20.84346284166959
$arr[14,0] = 10.18363778084792
$arr[14,1] = 5.476705462164446
$arr[14,2] = 11.99927614080916
$arr[14,3] = 0
$arr[14,4] = 29.38624282684169
$arr[14,5] = 27.10920186642424
$arr[14,6] = 13.85874358865237
$arr[14,7] = 19.55627602489837
$arr[14,8] = 11.42528758612231
$arr[14,9] = 9.723559724943986
$arr[14,10] = 0
$arr[14,11] = 0
$arr[14,12] = 17.66710028743574
$arr[14,13] = 20.87826590555955
$arr[15,0] = 10.05077209060246
$arr[15,1] = 5.390188157642988
$arr[15,2] = 11.94584725620297
$arr[15,3] = 0
$arr[15,4] = 29.38251270457535
$arr[15,5] = 27.119999706285
$arr[15,6] = 13.87405618702779
$arr[15,7] = 19.58698706249385
$arr[15,8] = 11.40975252464873
$arr[15,9] = 9.62850491430536
$arr[15,10] = 0
$arr[15,11] = 0
$arr[15,12] = 17.6863147327519
$arr[15,13] = 20.90072419159549
$arr[16,0] = 9.973609142508478
$arr[16,1] = 5.339705823235476
$arr[16,2] = 11.91523305014402
$arr[16,3] = 0
$arr[16,4] = 29.38105891517268
$arr[16,5] = 27.12702856753444
$arr[16,6] = 13.88305918841327
$arr[16,7] = 19.60498876923659
$arr[16,8] = 11.4010251328426
$arr[16,9] = 9.573477013365094
$arr[16,10] = 0
$arr[16,11] = 0
$arr[16,12] = 17.69751361547346
$arr[16,13] = 20.91404825302489
$arr[17,0] = 9.947357926828323
$arr[17,1] = 5.322490319762977
$arr[17,2] = 11.90488855249025
$arr[17,3] = 0
$arr[17,4] = 29.38068554964048
$arr[17,5] = 27.12954883347365
$arr[17,6] = 13.88614104657901
$arr[17,7] = 19.61114180107683
$arr[17,8] = 11.39810608072862
$arr[17,9] = 9.554786349677487
$arr[17,10] = 0
$arr[17,11] = 0
$arr[17,12] = 17.70133067904442
$arr[17,13] = 20.91862935133847
$arr[18,0] = 10.064993214793
$arr[18,1] = 5.399472712727194
$arr[18,2] = 11.9515229925237
$arr[18,3] = 0
$arr[18,4] = 29.38283821269431
$arr[18,5] = 27.11876554687174
$arr[18,6] = 13.87240589144297
$arr[18,7] = 19.58368288509925
$arr[18,8] = 11.41138477673519
$arr[18,9] = 9.638660813068798
$arr[18,10] = 0
$arr[18,11] = 0
$arr[18,12] = 17.68425408926576
$arr[18,13] = 20.89829137199625
$arr[19,0] = 10.45120408642595
$arr[19,1] = 5.64948225816953
$arr[19,2] = 12.10953514510156
$arr[19,3] = 0
$arr[19,4] = 29.39832917790197
$arr[19,5] = 27.09231902455194
$arr[19,6] = 13.82859448939132
$arr[19,7] = 19.49545805482082
$arr[19,8] = 11.45847076689022
$arr[19,9] = 9.91609934837545
$arr[19,10] = 0
$arr[19,11] = 0
$arr[19,12] = 17.62863644549525
$arr[19,13] = 20.83482019332972
$arr[20,0] = 10.6965367448226
$arr[20,1] = 5.806353433678748
$arr[20,2] = 12.21364084261683
$arr[20,3] = 0
$arr[20,4] = 29.41459362710742
$arr[20,5] = 27.08234261404838
$arr[20,6] = 13.80170254099406
$arr[20,7] = 19.44080562508609
$arr[20,8] = 11.49105489011301
$arr[20,9] = 10.09388266096328
$arr[20,10] = 0
$arr[20,11] = 0
$arr[20,12] = 17.59359266305292
$arr[20,13] = 20.79696844480422
$arr[21,0] = 10.56626686019709
$arr[21,1] = 5.723230224191737
$arr[21,2] = 12.15801104329625
$arr[21,3] = 0
$arr[21,4] = 29.40535225087404
$arr[21,5] = 27.08699767354739
$arr[21,6] = 13.81589596905039
$arr[21,7] = 19.46969936009344
$arr[21,8] = 11.47350052476083
$arr[21,9] = 9.999337598863121
$arr[21,10] = 0
$arr[21,11] = 0
$arr[21,12] = 17.61217699898174
$arr[21,13] = 20.81683802844968
$arr[22,0] = 10.05856625832246
$arr[22,1] = 5.395277476404478
$arr[22,2] = 11.94895666974322
$arr[22,3] = 0
$arr[22,4] = 29.38268889878126
$arr[22,5] = 27.11932095313173
$arr[22,6] = 13.87315136837463
$arr[22,7] = 19.58517562743833
$arr[22,8] = 11.41064619992953
$arr[22,9] = 9.634070509347579
$arr[22,10] = 0
$arr[22,11] = 0
$arr[22,12] = 17.68518523189407
$arr[22,13] = 20.89938996494617
$arr[23,0] = 9.482922039260245
$arr[23,1] = 5.014107913848855
$arr[23,2] = 11.72803727228428
$arr[23,3] = 0
$arr[23,4] = 29.38468936844314
$arr[23,5] = 27.18544707354917
$arr[23,6] = 13.94238915897161
$arr[23,7] = 19.72263932681463
$arr[23,8] = 11.35081842408121
$arr[23,9] = 9.226770622371271
$arr[23,10] = 0
$arr[23,11] = 0
$arr[23,12] = 17.76956944014304
$arr[23,13] = 20.96667396946135

$range = $ws.Range("B2:O25")
$range.Value2 = $arr